$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.22%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.73%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.924"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.83%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06546"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.41%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.246"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.32%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.366"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'12.88%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9145"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.33%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1573"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'4.03%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06530"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'27.53%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07655"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.57%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02980"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.84%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08986"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.04%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001595"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.34%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006540"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.58%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006108"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.41%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.486"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.57%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.397"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.56%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.242"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.84%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.66%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.00%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.15%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04480"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.60%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'10.10%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001185"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.47%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004334"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'12.42%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D28").Value = "'0.0001180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-1.72%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0001636"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'-15.73%"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Value = "'0.04148"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.52%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006980"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.77%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1416"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'20.61%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-1.55%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01249"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'8.95%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005539"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.81%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-6.93%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-7.56%"
$ws.Range("E47").Style = "Normal"
